$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header/description text as a new shared string, placed in D2.
$ws.Range("D2").Value = "File name\Extensions seperated by a '/' for every entry"

# Move the active selection to D3, matching the post-edit cursor position.
$ws.Range("D3").Select()
